# Actualización desde MV -datos-
# Adds a new "01-06-2021" monthly row (row 66) to Sheet1 and revises the
# previous last row's (row 65, "01-05-2021") values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revise the existing last data row (01-05-2021).
$ws.Range("B65").Value = -0.4
$ws.Range("C65").Value = -0.6
$ws.Range("D65").Value = -0.4

# Append the new row (01-06-2021). Force column A to be entered as literal
# text (matching the "mm-dd-yyyy"-looking strings already used throughout
# column A) instead of letting Excel auto-convert it to a date serial.
$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = "01-06-2021"
$ws.Range("A66").Style = "Normal"

$ws.Range("B66").Value = 0.3
$ws.Range("C66").Value = 0.2
$ws.Range("D66").Value = 0.4
